$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add "Errors" sheet right after "Classes"
$errors = $wb.Worksheets.Add($null, $ws1)
$errors.Name = "Errors"

# Add "Warnings" sheet right after "Errors"
$warnings = $wb.Worksheets.Add($null, $errors)
$warnings.Name = "Warnings"

$errors.Range("A1").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [001] ""Lunch""'',"'
$errors.Range("A2").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [102] ""PreK""'',"'
$errors.Range("A3").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [011] ""Kindergarten""'',"'
$errors.Range("A4").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [101] ""First Grade""'',"'
$errors.Range("A5").Formula = '="''Sheet ""Classes"" A subclass with the id ""8003"" was not found for Class [101] ""First Grade""'',"'
$errors.Range("A6").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [201] ""Second Grade""'',"'
$errors.Range("A7").Formula = '="''Sheet ""Classes"" A subclass with the id ""8003"" was not found for Class [201] ""Second Grade""'',"'
$errors.Range("A8").Formula = '="''Sheet ""Classes"" A subclass with the id ""8002"" was not found for Class [301] ""Third Grade""'',"'
$errors.Range("A9").Formula = '="''Sheet ""Classes"" A subclass with the id ""8003"" was not found for Class [301] ""Third Grade""'',"'

# Convert formulas to plain values (shared strings) via copy/paste-special
$errors.Range("A1:A9").Copy()
$errors.Range("A1:A9").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Remove tabSelected from Classes sheet, make Errors the active/selected tab
$ws1.Range("B2").Select()
$errors.Activate()
$errors.Range("B14").Select()
